$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E22").NumberFormat = "@"

$ws.Range("A22").Value = "edit1"
$ws.Range("B22").Value = "riya-morankar"
$ws.Range("C22").Value = "Merged"
$ws.Range("D22").Value = "cleared"
$ws.Range("E22").Value = "2025-06-30"
$ws.Range("F22").Value = "d9e149383f13845f5987656bf79fd690e8faab53"
